# Weekly update: insert two new price records (rows) into the Zanahoria
# (carrot) price sheet for Vega Monumental Concepcion.
#
# The existing data block occupies rows 180-283 (dated older entries).
# Two brand-new rows of data are inserted:
#   - one becomes the new row 180 (pushing the former row180-188 down to 181-189)
#   - one becomes the new row 190 (pushing everything from the former row189 on down by one more)
# giving a net growth of the sheet from 283 to 285 data/header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new row at position 180 ---------------------------------
$ws.Rows.Item(180).Insert()

$row = 180
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44846
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100114013
$ws.Cells.Item($row, 7).Value = "Zanahoria"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 270
$ws.Cells.Item($row, 11).Value = 21000
$ws.Cells.Item($row, 12).Value = 22000
$ws.Cells.Item($row, 13).Value = 21556
$ws.Cells.Item($row, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item($row, 15).Value = "Región de La Araucanía"
$ws.Cells.Item($row, 16).Value = 1078
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"

# --- Insert second new row at position 190 ---------------------------------
$ws.Rows.Item(190).Insert()

$row = 190
$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item($row, 3).Value = "Bíobío"
$ws.Cells.Item($row, 4).Value = 44845
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 100114013
$ws.Cells.Item($row, 7).Value = "Zanahoria"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 220
$ws.Cells.Item($row, 11).Value = 15000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 15545
$ws.Cells.Item($row, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item($row, 15).Value = "Región de La Araucanía"
$ws.Cells.Item($row, 16).Value = 777
$ws.Cells.Item($row, 17).Value = 20
$ws.Cells.Item($row, 18).Value = "Hortaliza"
